$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 9.718953170468685
$ws.Cells.Item(2, 3).Value = -8.452308629668783
$ws.Cells.Item(2, 4).Value = -0.3126370826896341
$ws.Cells.Item(2, 5).Value = 1.834933389109758
$ws.Cells.Item(2, 6).Value = -1.605096597413572
$ws.Cells.Item(2, 7).Value = -1.666270775827773
$ws.Cells.Item(2, 8).Value = 0.3751851613284818
$ws.Cells.Item(2, 9).Value = -0.2074093757051592
$ws.Cells.Item(2, 10).Value = -0.1117343348958014
$ws.Cells.Item(2, 11).Value = -0.3721533718945242

# Row 3
$ws.Cells.Item(3, 2).Value = -10.70211146928018
$ws.Cells.Item(3, 3).Value = -2.562439922301026
$ws.Cells.Item(3, 4).Value = -0.4148694505016339
$ws.Cells.Item(3, 5).Value = -3.854899437024964
$ws.Cells.Item(3, 6).Value = -3.916073615439165
$ws.Cells.Item(3, 7).Value = -1.87461767828291
$ws.Cells.Item(3, 8).Value = -2.457212215316551
$ws.Cells.Item(3, 9).Value = -2.361537174507194
$ws.Cells.Item(3, 10).Value = -2.621956211505916
$ws.Cells.Item(3, 11).Value = -1.786501215347425

# Row 4
$ws.Cells.Item(4, 2).Value = -10.03159025315832
$ws.Cells.Item(4, 3).Value = -7.884019781358927
$ws.Cells.Item(4, 4).Value = -11.32404976788226
$ws.Cells.Item(4, 5).Value = -11.38522394629646
$ws.Cells.Item(4, 6).Value = -9.343768009140204
$ws.Cells.Item(4, 7).Value = -9.926362546173845
$ws.Cells.Item(4, 8).Value = -9.830687505364487
$ws.Cells.Item(4, 9).Value = -10.09110654236321
$ws.Cells.Item(4, 10).Value = -9.255651546204719
$ws.Cells.Item(4, 11).Value = -9.35573090610734

# Row 5
$ws.Cells.Item(5, 2).Value = 10.28724201877854
$ws.Cells.Item(5, 3).Value = 6.847212032255212
$ws.Cells.Item(5, 4).Value = 6.786037853841011
$ws.Cells.Item(5, 5).Value = 8.827493790997265
$ws.Cells.Item(5, 6).Value = 8.244899253963624
$ws.Cells.Item(5, 7).Value = 8.340574294772981
$ws.Cells.Item(5, 8).Value = 8.080155257774258
$ws.Cells.Item(5, 9).Value = 8.91561025393275
$ws.Cells.Item(5, 10).Value = 8.815530894030129
$ws.Cells.Item(5, 11).Value = 8.461744387950432

# Row 6
$ws.Cells.Item(6, 2).Value = -1.292459514723937
$ws.Cells.Item(6, 3).Value = -1.353633693138139
$ws.Cells.Item(6, 4).Value = 0.6878222440181159
$ws.Cells.Item(6, 5).Value = 0.1052277069844749
$ws.Cells.Item(6, 6).Value = 0.2009027477938327
$ws.Cells.Item(6, 7).Value = -0.0595162892048901
$ws.Cells.Item(6, 8).Value = 0.7759387069536011
$ws.Cells.Item(6, 9).Value = 0.6758593470509796
$ws.Cells.Item(6, 10).Value = 0.3220728409712834
$ws.Cells.Item(6, 11).Value = 0.4702502314037287

# Row 7
$ws.Cells.Item(7, 2).Value = -3.501204164937531
$ws.Cells.Item(7, 3).Value = -1.459748227781277
$ws.Cells.Item(7, 4).Value = -2.042342764814918
$ws.Cells.Item(7, 5).Value = -1.94666772400556
$ws.Cells.Item(7, 6).Value = -2.207086761004283
$ws.Cells.Item(7, 7).Value = -1.371631764845791
$ws.Cells.Item(7, 8).Value = -1.471711124748413
$ws.Cells.Item(7, 9).Value = -1.825497630828109
$ws.Cells.Item(7, 10).Value = -1.677320240395664
$ws.Cells.Item(7, 11).Value = -1.519874251391509

# Row 8
$ws.Cells.Item(8, 2).Value = 1.980281758742053
$ws.Cells.Item(8, 3).Value = 1.397687221708412
$ws.Cells.Item(8, 4).Value = 1.49336226251777
$ws.Cells.Item(8, 5).Value = 1.232943225519047
$ws.Cells.Item(8, 6).Value = 2.068398221677539
$ws.Cells.Item(8, 7).Value = 1.968318861774917
$ws.Cells.Item(8, 8).Value = 1.614532355695221
$ws.Cells.Item(8, 9).Value = 1.762709746127666
$ws.Cells.Item(8, 10).Value = 1.920155735131821
$ws.Cells.Item(8, 11).Value = 1.41949627860575

# Row 9
$ws.Cells.Item(9, 2).Value = 1.458861400122613
$ws.Cells.Item(9, 3).Value = 1.554536440931971
$ws.Cells.Item(9, 4).Value = 1.294117403933248
$ws.Cells.Item(9, 5).Value = 2.12957240009174
$ws.Cells.Item(9, 6).Value = 2.029493040189118
$ws.Cells.Item(9, 7).Value = 1.675706534109422
$ws.Cells.Item(9, 8).Value = 1.823883924541867
$ws.Cells.Item(9, 9).Value = 1.981329913546022
$ws.Cells.Item(9, 10).Value = 1.480670457019951
$ws.Cells.Item(9, 11).Value = 1.761028337555421

# Row 10
$ws.Cells.Item(10, 2).Value = -0.4869194962242832
$ws.Cells.Item(10, 3).Value = -0.747338533223006
$ws.Cells.Item(10, 4).Value = 0.0881164629354852
$ws.Cells.Item(10, 5).Value = -0.01196289696713632
$ws.Cells.Item(10, 6).Value = -0.3657494030468326
$ws.Cells.Item(10, 7).Value = -0.2175720126143872
$ws.Cells.Item(10, 8).Value = -0.06012602361023223
$ws.Cells.Item(10, 9).Value = -0.560785480136303
$ws.Cells.Item(10, 10).Value = -0.2804275996008339
$ws.Cells.Item(10, 11).Value = -0.4524024574458455

# Row 11
$ws.Cells.Item(11, 2).Value = -0.164743996189365
$ws.Cells.Item(11, 3).Value = 0.6707109999691262
$ws.Cells.Item(11, 4).Value = 0.5706316400665047
$ws.Cells.Item(11, 5).Value = 0.2168451339868084
$ws.Cells.Item(11, 6).Value = 0.3650225244192538
$ws.Cells.Item(11, 7).Value = 0.5224685134234088
$ws.Cells.Item(11, 8).Value = 0.02180905689733798
$ws.Cells.Item(11, 9).Value = 0.3021669374328071
$ws.Cells.Item(11, 10).Value = 0.1301920795877955
$ws.Cells.Item(11, 11).Value = 0.436054619334127

# Row 12
$ws.Cells.Item(12, 2).Value = 0.5750359591597685
$ws.Cells.Item(12, 3).Value = 0.4749565992571469
$ws.Cells.Item(12, 4).Value = 0.1211700931774507
$ws.Cells.Item(12, 5).Value = 0.269347483609896
$ws.Cells.Item(12, 6).Value = 0.426793472614051
$ws.Cells.Item(12, 7).Value = -0.07386598391201982
$ws.Cells.Item(12, 8).Value = 0.2064918966234494
$ws.Cells.Item(12, 9).Value = 0.0345170387784377
$ws.Cells.Item(12, 10).Value = 0.3403795785247692
$ws.Cells.Item(12, 11).Value = -0.271950374162517

# Row 13
$ws.Cells.Item(13, 2).Value = 0.7353756362558697
$ws.Cells.Item(13, 3).Value = 0.3815891301761735
$ws.Cells.Item(13, 4).Value = 0.5297665206086188
$ws.Cells.Item(13, 5).Value = 0.6872125096127738
$ws.Cells.Item(13, 6).Value = 0.186553053086703
$ws.Cells.Item(13, 7).Value = 0.4669109336221722
$ws.Cells.Item(13, 8).Value = 0.2949360757771605
$ws.Cells.Item(13, 9).Value = 0.600798615523492
$ws.Cells.Item(13, 10).Value = -0.01153133716379418
$ws.Cells.Item(13, 11).Value = 0.6768400480353174

# Row 14
$ws.Cells.Item(14, 2).Value = -0.4538658659823178
$ws.Cells.Item(14, 3).Value = -0.3056884755498724
$ws.Cells.Item(14, 4).Value = -0.1482424865457174
$ws.Cells.Item(14, 5).Value = -0.6489019430717882
$ws.Cells.Item(14, 6).Value = -0.3685440625363191
$ws.Cells.Item(14, 7).Value = -0.5405189203813308
$ws.Cells.Item(14, 8).Value = -0.2346563806349992
$ws.Cells.Item(14, 9).Value = -0.8469863333222853
$ws.Cells.Item(14, 10).Value = -0.1586149481231739
$ws.Cells.Item(14, 11).Value = -0.4353061035472806

# Row 15
$ws.Cells.Item(15, 2).Value = -0.2056091156472509
$ws.Cells.Item(15, 3).Value = -0.04816312664309591
$ws.Cells.Item(15, 4).Value = -0.5488225831691667
$ws.Cells.Item(15, 5).Value = -0.2684647026336975
$ws.Cells.Item(15, 6).Value = -0.4404395604787092
$ws.Cells.Item(15, 7).Value = -0.1345770207323777
$ws.Cells.Item(15, 8).Value = -0.7469069734196638
$ws.Cells.Item(15, 9).Value = -0.05853558822055238
$ws.Cells.Item(15, 10).Value = -0.3352267436446591
$ws.Range("K15:K15").ClearContents()

# Row 16
$ws.Cells.Item(16, 2).Value = 0.3056233794366003
$ws.Cells.Item(16, 3).Value = -0.1950360770894705
$ws.Cells.Item(16, 4).Value = 0.08532180344599868
$ws.Cells.Item(16, 5).Value = -0.08665305439901295
$ws.Cells.Item(16, 6).Value = 0.2192094853473185
$ws.Cells.Item(16, 7).Value = -0.3931204673399676
$ws.Cells.Item(16, 8).Value = 0.2952509178591439
$ws.Cells.Item(16, 9).Value = 0.01855976243503714
$ws.Range("J16:K16").ClearContents()

# Row 17
$ws.Cells.Item(17, 2).Value = -0.3432134675219158
$ws.Cells.Item(17, 3).Value = -0.06285558698644665
$ws.Cells.Item(17, 4).Value = -0.2348304448314583
$ws.Cells.Item(17, 5).Value = 0.0710320949148732
$ws.Cells.Item(17, 6).Value = -0.541297857772413
$ws.Cells.Item(17, 7).Value = 0.1470735274266985
$ws.Cells.Item(17, 8).Value = -0.1296176279974082
$ws.Range("I17:K17").ClearContents()

# Row 18
$ws.Cells.Item(18, 2).Value = -0.2203015759906016
$ws.Cells.Item(18, 3).Value = -0.3922764338356133
$ws.Cells.Item(18, 4).Value = -0.0864138940892818
$ws.Cells.Item(18, 5).Value = -0.698743846776568
$ws.Cells.Item(18, 6).Value = -0.01037246157745647
$ws.Cells.Item(18, 7).Value = -0.2870636170015632
$ws.Range("H18:K18").ClearContents()

# Row 19
$ws.Cells.Item(19, 2).Value = 0.1083830226904575
$ws.Cells.Item(19, 3).Value = 0.414245562436789
$ws.Cells.Item(19, 4).Value = -0.1980843902504972
$ws.Cells.Item(19, 5).Value = 0.4902869949486143
$ws.Cells.Item(19, 6).Value = 0.2135958395245076
$ws.Range("G19:K19").ClearContents()

# Row 20
$ws.Cells.Item(20, 2).Value = 0.1338876819013198
$ws.Cells.Item(20, 3).Value = -0.4784422707859664
$ws.Cells.Item(20, 4).Value = 0.2099291144131452
$ws.Cells.Item(20, 5).Value = -0.06676204101096155
$ws.Range("F20:K20").ClearContents()

# Row 21
$ws.Cells.Item(21, 2).Value = -0.3064674129409547
$ws.Cells.Item(21, 3).Value = 0.3819039722581568
$ws.Cells.Item(21, 4).Value = 0.1052128168340501
$ws.Range("E21:K21").ClearContents()

# Row 22
$ws.Cells.Item(22, 2).Value = 0.07604143251182532
$ws.Cells.Item(22, 3).Value = -0.2006497229122814
$ws.Range("D22:K22").ClearContents()

# Row 23
$ws.Cells.Item(23, 2).Value = 0.4116802297750048
$ws.Range("C23:K23").ClearContents()

# Row 24
$ws.Range("B24:K24").ClearContents()
